$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. "batch_header" table (column A): insert a new "year" field
#        before "zoom id", pushing created_by/created_date down one row.
$ws.Range("A9").Value = "created_date"
$ws.Range("A8").Value = "created_by"
$ws.Range("A7").Value = "zoom id"
$ws.Range("A6").Value = "year"

# --- 2. Capitalise the "teacher" table header in column E
$ws.Range("E7").Value = "Teacher"

# --- 3. New "Attendance" table in columns L:M, mirroring the style of
#        the other little ER-diagram tables on the sheet.
$ws.Range("L1").Value = "Attendance_header"
$ws.Range("L1").Font.Bold = $true

$ws.Range("L2").Value = "att_id"
$ws.Range("M2").Value = "(P)"

$ws.Range("L3").Value = "bh_id"
$ws.Range("M3").Value = "(F),(I)"

$ws.Range("L4").Value = "att_pass"

$ws.Range("L5").Value = "created_by"

$ws.Range("L6").Value = "created_date"

$ws.Range("L9").Value = "Attendance_details"
$ws.Range("L9").Font.Bold = $true

$ws.Range("L10").Value = "att_id"
$ws.Range("M10").Value = "(F),(I)"

$ws.Range("L11").Value = "stu_id"
$ws.Range("M11").Value = "(F),(I)"

$ws.Range("L12").Value = "att_status"

# --- 4. Selection cursor moved back to A7 (matches the saved view state)
[void]$ws.Range("A7").Select()
